$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the wealth-tax model parameters (h, m, f)
$ws.Range("B2").Value = 1.36350750271454
$ws.Range("B3").Value = 1.02478093770152
$ws.Range("B4").Value = 231936.02942000001

# Update the "wealth" simulation input row
$ws.Range("B41").Value = 0.50105135999999995
$ws.Range("C41").Value = 0.65197760999999999
$ws.Range("D41").Value = 0.75971648000000003
$ws.Range("E41").Value = 0.84668591000000004
$ws.Range("F41").Value = 0.95077845000000005
$ws.Range("G41").Value = 1.2536148899999999
$ws.Range("H41").Value = 3.4556963500000002

# Keep the manual sim_conv literal in sync with f (A42)
$ws.Range("A42").Value = 231936.02942000001

# Move the active selection as recorded in the saved file
$ws.Range("B13").Select()
